# Append new daily cumulative-deaths records (rows 160-181) to Sheet1,
# mirroring the style of the existing data (column A uses the date style
# already applied to A2:A159 via the "s=1" cellXf).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @(44279, 9313, 1733, 11046),
    @(44280, 9373, 1742, 11115),
    @(44281, 9426, 1745, 11171),
    @(44282, 9496, 1752, 11248),
    @(44283, 9542, 1759, 11301),
    @(44284, 9624, 1782, 11406),
    @(44285, 9719, 1813, 11532),
    @(44286, 9790, 1837, 11627),
    @(44287, 9877, 1856, 11733),
    @(44288, 9948, 1861, 11809),
    @(44289, 10025, 1868, 11893),
    @(44290, 10094, 1877, 11971),
    @(44291, 10156, 1888, 12044),
    @(44292, 10243, 1894, 12137),
    @(44293, 10322, 1918, 12240),
    @(44294, 10411, 1949, 12360),
    @(44295, 10487, 1976, 12463),
    @(44296, 10552, 1990, 12542),
    @(44297, 10565, 2002, 12567),
    @(44298, 10716, 2043, 12759),
    @(44299, 10798, 2060, 12858),
    @(44300, 10877, 2090, 12967)
)

$startRow = 160
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Reset the view: scroll back to the top and select A1 (matches the
# simplified <sheetView .../> with no topLeftCell/selection override).
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
